# The commit swaps the presentation's applied theme color palette from the
# "Integral" design's "Red Violet" scheme over to the stock "Office Theme"
# palette (the same palette that ships as PowerPoint's default "Office"
# color scheme, which in this deck already lived - unused - as the
# notes-master's theme). Font scheme / format scheme are identical between
# the two themes, so the only observable content change is the 12 theme
# colors (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

function Set-ThemeColor {
    param($Scheme, [int]$Index, [int]$R, [int]$G, [int]$B)
    $Scheme.Item($Index).RGB = $R + ($G * 256) + ($B * 65536)
}

# Target palette = stock "Office" theme colors, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
Set-ThemeColor $cs 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $cs 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $cs 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $cs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $cs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $cs 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $cs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $cs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $cs 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $cs 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $cs 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $cs 12 0x95 0x4F 0x72   # folHlink
